# Test data for Greece Market
#
# Adds a new "Greece" worksheet at the end of the workbook by duplicating
# the existing "Croatia" sheet (same layout/styles) and updating the two
# market-specific cells to the Greece values. The new sheet becomes the
# active tab, matching how Excel leaves a freshly copied sheet selected.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate "Croatia" and place the copy immediately after it (i.e. at the
# very end of the tab strip).
$croatia.Copy($null, $croatia)

# The copy is now the last sheet and is the active sheet.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market name / ticket-reference cells for Greece.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3168"
